$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.440985666666667
$ws.Range("H2").Value = 4.322957000000001
$ws.Range("I2").Value = 0.1098365531732288
$ws.Range("J2").Value = 0.1230162332390494
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06706333333333332
$ws.Range("N2").Value = 0.20119
$ws.Range("O2").Value = 0.01564803973866319
$ws.Range("P2").Value = 0.01878890697104963
$ws.Range("Q2").Value = 0.09663730209222222
$ws.Range("R2").Value = 0.86973571883
$ws.Range("S2").Value = 0.001718726748812476
$ws.Range("T2").Value = 0.002311340562257442
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.440985666666667
$ws.Range("H3").Value = 4.322957000000001
$ws.Range("I3").Value = 0.1098365531732288
$ws.Range("J3").Value = 0.1230162332390494
$ws.Range("O3").Value = 0.2726197454399388
$ws.Range("P3").Value = 0.3273398534952746
$ws.Range("Q3").Value = 1.683612588948778
$ws.Range("R3").Value = 15.152513300539
$ws.Range("S3").Value = 0.02994361316608594
$ws.Range("T3").Value = 0.04026811576601095
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.440985666666667
$ws.Range("H4").Value = 4.322957000000001
$ws.Range("I4").Value = 0.1098365531732288
$ws.Range("J4").Value = 0.1230162332390494
$ws.Range("M4").Value = 0.72155
$ws.Range("N4").Value = 2.16465
$ws.Range("O4").Value = 0.1683608987539007
$ws.Range("P4").Value = 0.2021542197668005
$ws.Range("Q4").Value = 1.039743207783334
$ws.Range("R4").Value = 9.357688870050001
$ws.Range("S4").Value = 0.0184921808082754
$ws.Range("T4").Value = 0.02486825064909078
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.440985666666667
$ws.Range("H5").Value = 4.322957000000001
$ws.Range("I5").Value = 0.1098365531732288
$ws.Range("J5").Value = 0.1230162332390494
$ws.Range("M5").Value = 2.1492875
$ws.Range("N5").Value = 4.298575
$ws.Range("O5").Value = 0.501498129277977
$ws.Range("P5").Value = 0.4014390664699025
$ws.Range("Q5").Value = 3.097092481045834
$ws.Range("R5").Value = 18.582554886275
$ws.Range("S5").Value = 0.05508282594271529
$ws.Range("T5").Value = 0.04938352183212776
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.440985666666667
$ws.Range("H6").Value = 4.322957000000001
$ws.Range("I6").Value = 0.1098365531732288
$ws.Range("J6").Value = 0.1230162332390494
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1794573333333333
$ws.Range("N6").Value = 0.538372
$ws.Range("O6").Value = 0.04187318678952025
$ws.Range("P6").Value = 0.05027795329697268
$ws.Range("Q6").Value = 0.2585954451115556
$ws.Range("R6").Value = 2.327359006004
$ws.Range("S6").Value = 0.004599206507339683
$ws.Range("T6").Value = 0.006185004429562423
$ws.Range("I7").Value = 0.5687502547919595
$ws.Range("J7").Value = 0.6369966279614609
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06706333333333332
$ws.Range("N7").Value = 0.20119
$ws.Range("O7").Value = 0.01564803973866319
$ws.Range("P7").Value = 0.01878890697104963
$ws.Range("Q7").Value = 0.5004025399511109
$ws.Range("R7").Value = 4.503622859559999
$ws.Range("S7").Value = 0.008899826588359395
$ws.Range("T7").Value = 0.0119684703836402
$ws.Range("I8").Value = 0.5687502547919595
$ws.Range("J8").Value = 0.6369966279614609
$ws.Range("O8").Value = 0.2726197454399388
$ws.Range("P8").Value = 0.3273398534952746
$ws.Range("S8").Value = 0.1550525496802843
$ws.Range("T8").Value = 0.2085143828738886
$ws.Range("I9").Value = 0.5687502547919595
$ws.Range("J9").Value = 0.6369966279614609
$ws.Range("M9").Value = 0.72155
$ws.Range("N9").Value = 2.16465
$ws.Range("O9").Value = 0.1683608987539007
$ws.Range("P9").Value = 0.2021542197668005
$ws.Range("Q9").Value = 5.383947304066666
$ws.Range("R9").Value = 48.45552573659999
$ws.Range("S9").Value = 0.09575530406328431
$ws.Range("T9").Value = 0.128771556319632
$ws.Range("I10").Value = 0.5687502547919595
$ws.Range("J10").Value = 0.6369966279614609
$ws.Range("M10").Value = 2.1492875
$ws.Range("N10").Value = 4.298575
$ws.Range("O10").Value = 0.501498129277977
$ws.Range("P10").Value = 0.4014390664699025
$ws.Range("Q10").Value = 16.03721244721666
$ws.Range("R10").Value = 96.22327468329998
$ws.Range("S10").Value = 0.2852271888045405
$ws.Range("T10").Value = 0.2557153316733247
$ws.Range("I11").Value = 0.5687502547919595
$ws.Range("J11").Value = 0.6369966279614609
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1794573333333333
$ws.Range("N11").Value = 0.538372
$ws.Range("O11").Value = 0.04187318678952025
$ws.Range("P11").Value = 0.05027795329697268
$ws.Range("Q11").Value = 1.339046255969778
$ws.Range("R11").Value = 12.051416303728
$ws.Range("S11").Value = 0.02381538565549096
$ws.Range("T11").Value = 0.03202688671097541
$ws.Range("G12").Value = 4.2167365
$ws.Range("H12").Value = 8.433472999999999
$ws.Range("I12").Value = 0.3214131920348118
$ws.Range("J12").Value = 0.2399871387994896
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.06706333333333332
$ws.Range("N12").Value = 0.20119
$ws.Range("O12").Value = 0.01564803973866319
$ws.Range("P12").Value = 0.01878890697104963
$ws.Range("Q12").Value = 0.2827884054783333
$ws.Range("R12").Value = 1.69673043287
$ws.Range("S12").Value = 0.005029486401491317
$ws.Range("T12").Value = 0.004509096025151986
$ws.Range("G13").Value = 4.2167365
$ws.Range("H13").Value = 8.433472999999999
$ws.Range("I13").Value = 0.3214131920348118
$ws.Range("J13").Value = 0.2399871387994896
$ws.Range("O13").Value = 0.2726197454399388
$ws.Range("P13").Value = 0.3273398534952746
$ws.Range("Q13").Value = 4.926732319345167
$ws.Range("R13").Value = 29.560393916071
$ws.Range("S13").Value = 0.08762358259356856
$ws.Range("T13").Value = 0.07855735485537506
$ws.Range("G14").Value = 4.2167365
$ws.Range("H14").Value = 8.433472999999999
$ws.Range("I14").Value = 0.3214131920348118
$ws.Range("J14").Value = 0.2399871387994896
$ws.Range("M14").Value = 0.72155
$ws.Range("N14").Value = 2.16465
$ws.Range("O14").Value = 0.1683608987539007
$ws.Range("P14").Value = 0.2021542197668005
$ws.Range("Q14").Value = 3.042586221575
$ws.Range("R14").Value = 18.25551732945
$ws.Range("S14").Value = 0.05411341388234098
$ws.Range("T14").Value = 0.04851441279807768
$ws.Range("G15").Value = 4.2167365
$ws.Range("H15").Value = 8.433472999999999
$ws.Range("I15").Value = 0.3214131920348118
$ws.Range("J15").Value = 0.2399871387994896
$ws.Range("M15").Value = 2.1492875
$ws.Range("N15").Value = 4.298575
$ws.Range("O15").Value = 0.501498129277977
$ws.Range("P15").Value = 0.4014390664699025
$ws.Range("Q15").Value = 9.062979050243749
$ws.Range("R15").Value = 36.251916200975
$ws.Range("S15").Value = 0.1611881145307213
$ws.Range("T15").Value = 0.09634021296445001
$ws.Range("G16").Value = 4.2167365
$ws.Range("H16").Value = 8.433472999999999
$ws.Range("I16").Value = 0.3214131920348118
$ws.Range("J16").Value = 0.2399871387994896
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1794573333333333
$ws.Range("N16").Value = 0.538372
$ws.Range("O16").Value = 0.04187318678952025
$ws.Range("P16").Value = 0.05027795329697268
$ws.Range("Q16").Value = 0.7567242876593333
$ws.Range("R16").Value = 4.540345725956
$ws.Range("S16").Value = 0.01345859462668962
